$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-155). All of these were bumped from 45177 (2023-09-08) to
# 45178 (2023-09-09).
$range = $ws.Range("C2:C155")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45177) {
        $cell.Value = 45178
    }
}
